$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First "Vilhelm :" in the chat body -> bold "Vilhelm", keep " :" plain
# ------------------------------------------------------------------
$rng = $d.Range($d.Paragraphs(2).Range.Start, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Text = "Vilhelm"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchWholeWord = $true
if ($rng.Find.Execute()) {
    $rng.Font.Bold = 1
}

# ------------------------------------------------------------------
# 2) First "Johnny:" in the chat body (after the <br/>) -> bold "Johnny"
# ------------------------------------------------------------------
$rng = $d.Range($d.Paragraphs(2).Range.Start, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Text = "Johnny"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchWholeWord = $true
if ($rng.Find.Execute()) {
    $rng.Font.Bold = 1
}

# ------------------------------------------------------------------
# 3) "Vilhelm:" in the third paragraph -> bold "Vilhelm"
# ------------------------------------------------------------------
$rng = $d.Range($d.Paragraphs(3).Range.Start, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Text = "Vilhelm"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchWholeWord = $true
if ($rng.Find.Execute()) {
    $rng.Font.Bold = 1
}

# ------------------------------------------------------------------
# 4) Fourth (last, bookmark-only) paragraph gets Johnny's long reply.
#    "Johnny" is bold, the rest of the reply is plain text. The
#    _GoBack bookmark originally sits alone in that empty paragraph;
#    after the edit it sits in the middle of the new text (right
#    after "...alla sakerhetsprofil" and before "ter sa ar det; "),
#    so we insert the "before" text first (which pushes the bookmark
#    along with it), then relocate the bookmark, then append the
#    "after" text.
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range

$boldRun = $p4.Duplicate
$boldRun.Collapse(0)
$boldRun.InsertAfter("Johnny")
$boldRun.Font.Bold = 1

$restText = @"
: ” Vi har en ansvarig på varje plats, som får assistera oss vid utskjut av konfig till brandväggarna. Då har vi förbätt en konfiguration som läggs på ett USB minne och startar hårdvaran från det. Då kopplar den upp alla inställningar som vi på förhand har konfat, så kommer tunnel upp och vi driftsätter miljön. Så skickar ner en tom brandvägg, så används USB för att initisalt kofa upp WAN och en tunnel till deployment. Så blåser vi in en ny konfiguration som slår mot tunnel mot produktion och alla säkerhetsinställningar fixade. När det kommer till kontroll används; Statiska portar, applications skydd, webfilter, AV, IPS och DNS filtering i massa olika kombinationer för att uppnå säkerhetsnivåerna. Det skilljer sig tillexempel på mellan medarbetare och gäster, vi skyddar båda men på lite olika nivåer. Annars om man tar en best-practice på alla säkerhetsprofil
"@
$restRun = $d.Paragraphs(4).Range
$restRun.Collapse(0)
$restRun.InsertAfter($restText)
$restRun.Font.Bold = 0

# Relocate the hidden _GoBack bookmark to sit right before "ter sa ar det; "
$afterText = @"
ter så är det; 
"@
$splitPoint = $d.Paragraphs(4).Range
$splitPoint.Collapse(0)
$splitPoint.MoveStart(1, -1 * $afterText.Length)
$d.Bookmarks.Add("_GoBack", $splitPoint)

$tailRun = $d.Paragraphs(4).Range
$tailRun.Collapse(0)
$tailRun.InsertAfter($afterText)

# ------------------------------------------------------------------
# 5) Two more new trailing paragraphs.
# ------------------------------------------------------------------
$p4end = $d.Paragraphs(4).Range
$p4end.Collapse(0)
$p4end.InsertParagraphAfter()

$p5 = $d.Paragraphs(5).Range
$p5text = @"
Köra proxy för att kunna göra MIM för att kunna se SSL traifken, noga med rätt certificat och kunna lägga undantag för HTTPS eller SSH.
"@
$p5.InsertAfter($p5text)
$p5.Collapse(0)

$p5end = $d.Paragraphs(5).Range
$p5end.Collapse(0)
$p5end.InsertParagraphAfter()

$p6 = $d.Paragraphs(6).Range
$p6text = @"
Allt okänd skall blockeras, så att man går approch att allt som företaget anser som viktigt eller företagskritiskt så är det tillåtet från start. Sen kan man man få önskemål om att lägga undantag, i detta fall så är OFFICE 365 tillåtet med alla funktioner, men facebook är det begränsat med att chatt är tillåtet men inga APP, inga spel, inga video, inga file upload eller download.”
"@
$p6.InsertAfter($p6text)
